# [Fonds de solidarite] Add 2022-06-23 data
# Update nombre_aides (column C) and montant_total (column E) for the rows
# whose new figures were published with the 2022-06-23 data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 8;  C = 1050;   E = 91353904 },
    @{ Row = 63; C = 14364;  E = 36198911 },
    @{ Row = 64; C = 5218;   E = 20473121 },
    @{ Row = 65; C = 2018;   E = 13672438 },
    @{ Row = 71; C = 2601;   E = 5956007 },
    @{ Row = 72; C = 6279;   E = 15206435 },
    @{ Row = 91; C = 151171; E = 482664512 },
    @{ Row = 92; C = 409248; E = 1596880544 },
    @{ Row = 93; C = 209643; E = 1309884004 },
    @{ Row = 96; C = 17315;  E = 796624296 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
